# Refresh the crypto "Price" (D) and "Volume(1h)" (E) columns with the
# latest scraped snapshot values. Values are stored as literal text
# (e.g. "309.53", "-3.80%") in the source sheet, so we force the Text
# number format before assigning and then drop back to the "Normal"
# style afterwards -- this stops Excel from auto-coercing the strings
# into numbers/percentages (and from permanently leaving a custom
# number-format style applied to the cell).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "309.53"
Set-TextValue $ws.Range("E2") "-3.80%"
Set-TextValue $ws.Range("D3") "49.68"
Set-TextValue $ws.Range("E3") "2.03%"
Set-TextValue $ws.Range("D4") "5.160"
Set-TextValue $ws.Range("E4") "-2.08%"
Set-TextValue $ws.Range("D5") "0.07764"
Set-TextValue $ws.Range("E5") "-4.33%"
Set-TextValue $ws.Range("D6") "4.498"
Set-TextValue $ws.Range("E6") "-1.76%"
Set-TextValue $ws.Range("D7") "1.343"
Set-TextValue $ws.Range("E7") "11.53%"
Set-TextValue $ws.Range("D8") "1.559"
Set-TextValue $ws.Range("E8") "-5.30%"
Set-TextValue $ws.Range("D9") "0.1210"
Set-TextValue $ws.Range("E9") "-6.72%"
Set-TextValue $ws.Range("D10") "0.1982"
Set-TextValue $ws.Range("E10") "1.76%"
Set-TextValue $ws.Range("D11") "0.04785"
Set-TextValue $ws.Range("E11") "3.52%"
Set-TextValue $ws.Range("D12") "0.09376"
Set-TextValue $ws.Range("E12") "-0.36%"
Set-TextValue $ws.Range("E13") "-0.64%"
Set-TextValue $ws.Range("D14") "0.001267"
Set-TextValue $ws.Range("E14") "-4.10%"
Set-TextValue $ws.Range("D15") "0.005794"
Set-TextValue $ws.Range("E15") "-2.36%"
Set-TextValue $ws.Range("D16") "0.007491"
Set-TextValue $ws.Range("E16") "2,016.60%"
Set-TextValue $ws.Range("D17") "3.334"
Set-TextValue $ws.Range("D18") "2.433"
Set-TextValue $ws.Range("E18") "0.21%"
Set-TextValue $ws.Range("E19") "2.05%"
Set-TextValue $ws.Range("D20") "7.997"
Set-TextValue $ws.Range("E20") "-1.01%"
Set-TextValue $ws.Range("E21") "-2.28%"
Set-TextValue $ws.Range("D22") "0.3095"
Set-TextValue $ws.Range("E22") "-0.95%"
Set-TextValue $ws.Range("D23") "0.04167"
Set-TextValue $ws.Range("E23") "0.17%"
Set-TextValue $ws.Range("D24") "0.001273"
Set-TextValue $ws.Range("E24") "-2.52%"
Set-TextValue $ws.Range("D25") "0.003940"
Set-TextValue $ws.Range("E25") "-7.40%"
Set-TextValue $ws.Range("E26") "-0.01%"
Set-TextValue $ws.Range("D38") "0.02604"
Set-TextValue $ws.Range("E38") "-4.23%"
Set-TextValue $ws.Range("D39") "0.05980"
Set-TextValue $ws.Range("E39") "3.60%"
Set-TextValue $ws.Range("E40") "74.57%"
Set-TextValue $ws.Range("D41") "0.007934"
Set-TextValue $ws.Range("E41") "3.29%"
Set-TextValue $ws.Range("D42") "0.1423"
Set-TextValue $ws.Range("E42") "-1.29%"
Set-TextValue $ws.Range("D43") "0.008419"
Set-TextValue $ws.Range("E43") "9.43%"
Set-TextValue $ws.Range("D44") "0.008340"
Set-TextValue $ws.Range("E44") "2.95%"
Set-TextValue $ws.Range("D45") "0.3368"
Set-TextValue $ws.Range("E45") "5.45%"
Set-TextValue $ws.Range("D46") "0.00007604"
Set-TextValue $ws.Range("E46") "8.51%"
Set-TextValue $ws.Range("E47") "-0.02%"
Set-TextValue $ws.Range("E48") "-14.98%"
Set-TextValue $ws.Range("D49") "0.002619"
Set-TextValue $ws.Range("E49") "-34.51%"
Set-TextValue $ws.Range("D50") "0.00002100"
Set-TextValue $ws.Range("E50") "-0.02%"
Set-TextValue $ws.Range("D51") "0.0002000"
Set-TextValue $ws.Range("E51") "-0.02%"
